# ratio template: change label workerstation name
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the worker-station labels to their expanded/renamed forms.
$ws.Range("B5").Value = "Moldeado Cerámico"
$ws.Range("D5").Value = "Pintando Cerámico"
$ws.Range("F5").Value = "Horneado Cerámico"
$ws.Range("J5").Value = "Moldeado Retablo"
$ws.Range("L5").Value = "Pintado Retablo"

# Match the saved selection state in the target workbook.
$ws.Range("L6").Select()
